$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Clear the value in B3 (was "TestCase2") to simulate a null/blank value
$ws.Range("B3").ClearContents()

# Update the last active selection to B3
$ws.Range("B3").Select()
